$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update headers: C1 "Department" -> "Employee Id"; A1 "User Name (Employee Id)" -> "User Name"; add D1 "Role"
$ws.Range("C1").Value = "Employee Id"
$ws.Range("A1").Value = "User Name"
$ws.Range("D1").Value = "Role"

# Add a basic admin account row
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "tmml"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Admin"

# Set column C best-fit width to match its content (closest achievable to 11.1640625)
$ws.Columns("C").ColumnWidth = 10.3

# Select D2 to match final selection state
$ws.Range("D2").Select()
